# Moving from POI 3.17.0 to 4.0.1 - the embedded Java stack trace in the
# "setConserveRatioTrueInvalidImageFormat" expected-generation document
# shifts by a handful of source lines and the JUnit launcher tail changes
# from the Maven/Tycho/Equinox launcher chain to the Eclipse JDT JUnit
# runner chain.

$d = $word.ActiveDocument
$tab = [char]9

# --- Individual line-number shifts inside the stack trace ---
$lineEdits = @(
    @("org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:163)",
      "org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:162)"),
    @("org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:136)",
      "org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135)"),
    @("org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:168)",
      "org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)"),
    @("org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callOrApply(EvaluationServices.java:204)",
      "org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callOrApply(EvaluationServices.java:208)"),
    @("org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)",
      "org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)"),
    @("org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)",
      "org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:462)"),
    @("org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:389)",
      "org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:372)"),
    @("sun.reflect.GeneratedMethodAccessor74.invoke(Unknown Source)",
      "sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)")
)

foreach ($edit in $lineEdits) {
    $d.Content.Find.Execute($edit[0], $true, $false, $false, $false, $false, $true, 1, $false, $edit[1], 2) | Out-Null
}

# --- Replace the Maven/Tycho/Equinox launcher tail with the Eclipse JDT JUnit runner tail ---
$oldTailLines = @(
    "${tab}at org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)",
    "${tab}at org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)",
    "${tab}at org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)",
    "${tab}at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)",
    "${tab}at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)",
    "${tab}at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)",
    "${tab}at java.lang.reflect.Method.invoke(Method.java:498)",
    "${tab}at org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)",
    "${tab}at org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:156)",
    "${tab}at org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)",
    "${tab}at org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)",
    "${tab}at org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)",
    "${tab}at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)",
    "${tab}at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)",
    "${tab}at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)",
    "${tab}at java.lang.reflect.Method.invoke(Method.java:498)",
    "${tab}at org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)",
    "${tab}at org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)",
    "${tab}at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)",
    "${tab}at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)",
    "${tab}at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)",
    "${tab}at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)",
    "${tab}at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)",
    "${tab}at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)",
    "${tab}at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)",
    "${tab}at java.lang.reflect.Method.invoke(Method.java:498)",
    "${tab}at org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)",
    "${tab}at org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)",
    "${tab}at org.eclipse.equinox.launcher.Main.run(Main.java:1498)",
    "${tab}at org.eclipse.equinox.launcher.Main.main(Main.java:1471)"
)

$newTailLines = @(
    "${tab}at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)",
    "${tab}at org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)",
    "${tab}at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)",
    "${tab}at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)",
    "${tab}at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)",
    "${tab}at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)"
)

$oldTail = $oldTailLines -join "`n"
$newTail = $newTailLines -join "`n"

$d.Content.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2) | Out-Null
